$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.133.21'
$ws.Range('E2').Value = '  +1.26%  '

$ws.Range('D3').Value = '2.058.53'

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = '249.16'
$ws.Range('E5').Value = '  -1.92%  '

$ws.Range('D6').Value = '0.657'
$ws.Range('E6').Value = '  -1.16%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').Value = '55.66'
$ws.Range('E8').Value = '  +16.61%  '

$ws.Range('D9').Value = '61.52'
$ws.Range('E9').Value = '  +1.83%  '

$ws.Range('D10').Value = '0.381'
$ws.Range('E10').Value = '  +1.91%  '

$ws.Range('D11').Value = '0.0797'
$ws.Range('E11').Value = '  +6.72%  '

$ws.Range('D13').Value = '15.18'
$ws.Range('E13').Value = '  +6.03%  '

$ws.Range('D14').Value = '2.355.77'
$ws.Range('E14').Value = '  -2.64%  '

$ws.Range('D15').Value = '0.819'
$ws.Range('E15').Value = '  -1.76%  '

$ws.Range('D16').Value = '5.26'
$ws.Range('E16').Value = '  +2.63%  '

$ws.Range('D17').Value = '2.059.74'
$ws.Range('E17').Value = '  -2.37%  '

$ws.Range('D18').Value = '37.124.30'
$ws.Range('E18').Value = '  +1.17%  '

$ws.Range('D19').Value = '0.0₃0912'
$ws.Range('E19').Value = '  +9.15%  '

$ws.Range('D20').Value = '72.46'

$ws.Range('D21').Value = '14.27'
$ws.Range('E21').Value = '  +7.63%  '

$ws.Range('D22').Value = '5.37'
$ws.Range('E22').Value = '  +3.21%  '

$ws.Range('D23').Value = '237.28'
$ws.Range('E23').Value = '  -1.43%  '

$ws.Range('E24').Value = '  +0.01%  '

$ws.Range('E25').Value = '  -2.32%  '

$ws.Range('D26').Value = '170.10'
$ws.Range('E26').Value = '  -1.58%  '

$ws.Range('D27').Value = '9.06'
$ws.Range('E27').Value = '  -1.43%  '

$ws.Range('D28').Value = '20.21'
$ws.Range('E28').Value = '  -6.53%  '

$ws.Range('E29').Value = '  -1.74%  '

$ws.Range('E30').Value = '  -0.23%  '

$ws.Range('D31').Value = '4.56'
$ws.Range('E31').Value = '  +1.26%  '

$ws.Range('E32').Value = '  +9.64%  '

$ws.Range('D33').Value = '0.0624'
$ws.Range('E33').Value = '  +3.52%  '

$ws.Range('D34').Value = '4.34'
$ws.Range('E34').Value = '  +5.85%  '

$ws.Range('E35').Value = '  +0.07%  '

$ws.Range('D36').Value = '0.0864'
$ws.Range('E36').Value = '  -5.96%  '

$ws.Range('D37').Value = '2.28'
$ws.Range('E37').Value = '  -3.47%  '

$ws.Range('D38').Value = '1.77'
$ws.Range('E38').Value = '  -6.33%  '

$ws.Range('D39').Value = '1.36'
$ws.Range('E39').Value = '  +2.29%  '

$ws.Range('E40').Value = '  +22.13%  '

$ws.Range('D41').Value = '18.08'
$ws.Range('E41').Value = '  +12.89%  '

$ws.Range('D42').Value = '0.0224'
$ws.Range('E42').Value = '  -0.56%  '

$ws.Range('E43').Value = '  -3.97%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '96.66'
$ws.Range('E44').Value = '  -2.14%  '

$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').Value = '4.36'
$ws.Range('E45').Value = '  +46.83%  '

$ws.Range('B46').Value = 'HuobiToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D46').Value = '2.79'
$ws.Range('E46').Value = '  +0.11%  '

$ws.Range('B47').Value = 'Gas'
$ws.Range('C47').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D47').Value = '14.23'
$ws.Range('E47').Value = '  -52.19%  '

$ws.Range('E48').Value = '  +5.98%  '

$ws.Range('D49').Value = '1.298.82'
$ws.Range('E49').Value = '  -3.64%  '

$ws.Range('E50').Value = '  +2.68%  '

$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '4.04'
$ws.Range('E51').Value = '  +4.40%  '
